$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from E1 onto F1, then set the header text
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("F1").Value = "CHMO_DEF"

# CHMO_DEF values for rows 2-70 (aligned with CHMO_IRI in column D)
$chmoDef = @(
    '[''A generically dependent continuant that is about some thing. [IAO]'']',
    '[''Any constitutionally or isotopically distinct atom, molecule, ion, ion pair, radical, radical ion, complex, conformer etc., identifiable as a separately distinguishable entity. [IUPAC]'']',
    '[]',
    '[''A directive information entity that describes an intended process endpoint. When part of a plan specification the concretization is realized in a planned process in which the bearer tries to effect the world so that the process endpoint is achieved. [IAO]'']',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[''B is a role means: b is a realizable entity and b exists because there is some single bearer that is in some special physical, social, or institutional set of circumstances in which this bearer does not have to be and b is not such that, if it ceases to exist, then the physical make-up of the bearer is thereby changed. [BFO]'']',
    '[''A quality is a specifically dependent continuant that, in contrast to roles and dispositions, does not require any further process in order to be realized. [BFO]'']',
    '[]',
    '[]',
    '[''A material entity is an independent continuant that has some portion of matter as proper or improper continuant part. [BFO]'']',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[]',
    '[''b is an independent continuant if b is a continuant which is such that there is no c and no t such that b s-depends_on c at t. [BFO]'']',
    '[''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type.´[BFO]'', ''To say that b is a realizable entity is to say that b is a specifically dependent continuant that inheres in some independent continuant which is not a spatial region and is of a type instances of which are realized in processes of a correlated type. (axiom label in BFO2 Reference: [058-002])'']',
    '[''b is a generically dependent continuant if b is a continuant that generically depends on one or more other entities. [BFO]'']',
    '[''A continuant is an entity that persists, endures, or continues to exist through time while maintaining its identity. [BFO]'']',
    '[''b is a specifically dependent continuant if b is a continuant and there is some independent continuant c which is not a spatial region and which is such that b specifically depends on c at every time t during the course of b’s existence. [BFO]'']'
)

for ($i = 0; $i -lt $chmoDef.Length; $i++) {
    $ws.Cells.Item($i + 2, 6).Value = $chmoDef[$i]
}